$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "duplicate_image_filename" column (E) didn't have values for the main
# stimuli rows (2-21) -- fill them in with "NA".
$ws.Range("E2:E21").Value = "NA"
